$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new review as row 4 ---
$ws.Range("A4").Value = "com.singleton.strechy"
$ws.Range("B4").Value = "stretchy"
$ws.Range("C4").Value = "gazittalia1@gmail.com"
$ws.Range("D4").Value = "hermanliran@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:59"
$ws.Range("F4").Value = "genuine free offline car game"

# Turn the new email cells into live mailto: hyperlinks, same as rows 2 & 3
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:hermanliran@gmail.com", "", "", "hermanliran@gmail.com")

# Re-copy row 3's cell formatting onto row 4 so every column keeps the same
# look (font/alignment/etc.) as the rest of the table - a formats-only paste
# leaves the values and the hyperlinks just added untouched.
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)

# Adding the hyperlinks registered a built-in "Hyperlink" named style that is
# no longer used now that the original formatting was restored; drop it.
$wb.Styles("Hyperlink").Delete()

# Match the workbook's last saved selection (bottom-right filled cell)
$ws.Range("F4").Select()
